$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values replacing the old Strike# values in column G
$kValues = @{
    2  = 0
    3  = 0
    4  = 3
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 2
    11 = 0
    12 = 1
    13 = 2
    14 = 1
    17 = 1
    18 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
